$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avg Demand Scenario")

# Update Production capacity (O) and Facility cost (P) values for port scenario rows
$ws.Range("O11").Value = 250
$ws.Range("P11").Value = 250
$ws.Range("O12").Value = 275
$ws.Range("P12").Value = 250
$ws.Range("O13").Value = 275
$ws.Range("P13").Value = 250
$ws.Range("O14").Value = 275
$ws.Range("P14").Value = 250
$ws.Range("O15").Value = 275
$ws.Range("P15").Value = 250
$ws.Range("O16").Value = 100
$ws.Range("P16").Value = 250
$ws.Range("O17").Value = 100
$ws.Range("P17").Value = 250
$ws.Range("O18").Value = 100
$ws.Range("P18").Value = 250
$ws.Range("O19").Value = 100
$ws.Range("P19").Value = 250
$ws.Range("O20").Value = 150
$ws.Range("P20").Value = 300
$ws.Range("O21").Value = 100
$ws.Range("P21").Value = 250
$ws.Range("O22").Value = 100
$ws.Range("P22").Value = 250
$ws.Range("O23").Value = 100
$ws.Range("P23").Value = 250
$ws.Range("O24").Value = 100
$ws.Range("P24").Value = 400
$ws.Range("O26").Value = 50
$ws.Range("P26").Value = 10
$ws.Range("O27").Value = 50
$ws.Range("P27").Value = 10
$ws.Range("O28").Value = 100
$ws.Range("P28").Value = 200
$ws.Range("O30").Value = 1000
$ws.Range("P30").Value = 200
$ws.Range("O31").Value = 850
$ws.Range("P31").Value = 200
$ws.Range("O32").Value = 200
$ws.Range("P32").Value = 200
$ws.Range("O33").Value = 200
$ws.Range("P33").Value = 200
$ws.Range("O34").Value = 500
$ws.Range("P34").Value = 200
$ws.Range("O35").Value = 500
$ws.Range("P35").Value = 200
$ws.Range("O36").Value = 60000
$ws.Range("P36").Value = 1000
$ws.Range("O37").Value = 1000000
$ws.Range("P37").Value = 2000
$ws.Range("O38").Value = 700
$ws.Range("P38").Value = 100
$ws.Range("O39").Value = 5000
$ws.Range("P39").Value = 100
$ws.Range("O40").Value = 50
$ws.Range("P40").Value = 100
$ws.Range("O41").Value = 50
$ws.Range("P41").Value = 100
$ws.Range("O42").Value = 50
$ws.Range("P42").Value = 100
$ws.Range("O43").Value = 50
$ws.Range("P43").Value = 100
$ws.Range("O44").Value = 50
$ws.Range("P44").Value = 100
$ws.Range("O45").Value = 50
$ws.Range("P45").Value = 100
$ws.Range("O46").Value = 2000
$ws.Range("P46").Value = 50
$ws.Range("O47").Value = 2000
$ws.Range("P47").Value = 50

# Update the active window selection / scroll position to match the saved view
$ws.Activate()
$ws.Range("O11:P47").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 10

